$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "C2"  = 2.111455329893493
    "C3"  = 0.3920205665085499
    "C4"  = -0.014504764799976004
    "C5"  = 0.034827692605432635
    "C6"  = 0.06803680322205402
    "C7"  = 0.06792839823741074
    "C8"  = 0.022551890854185363
    "C9"  = -0.15717406572951634
    "C10" = 0.13454098486900007
    "C11" = -0.11289894807040429
    "C12" = -0.16436059978959433
    "C13" = 0.17943849699314326
    "C14" = 0.07835671027811729
    "C15" = -0.007718421908448897
    "C16" = -0.08035286371511285
    "C17" = 0.1051118163046719
    "C18" = -0.010810943718758377
}

foreach ($cellRef in $values.Keys) {
    $ws.Range($cellRef).Value = $values[$cellRef]
}
